# Generate Report for handoff
# Updates the localization-status workbook: a new source file
# (41f7d767-...md) replaces the old one (7e073687-...md), a brand-new
# file (ffff473fb4a2-...md) is added, and the zh-cn / de-de detail
# sheets record "Ready for handoff" / "Include" status together with
# the handoff package (.xlf) file + timestamp for both new rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Rebuild hyperlinks from scratch (clears the whole sheet collection).
$ws1.Range("A1").Hyperlinks.Delete()

$ws1.Range("A2").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/e2e/41f7d767-a907-42ee-b2c9-5f80e94a7af9.md", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/e2e/ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md", "", "", "ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f3066fa0cb0df4588b1a40478de660d455eaaa4d/.localization-config", "", "", ".localization-config") | Out-Null

# Restore the workbook's custom hyperlink look (underlined + blue) on
# the cells that carry a link, since Hyperlinks.Add applies its own
# default style.
$ws1.Range("A2").Font.Underline = 2
$ws1.Range("A2").Font.Color = 15570276
$ws1.Range("A3").Font.Underline = 2
$ws1.Range("A3").Font.Color = 15570276
$ws1.Range("A4").Font.Underline = 2
$ws1.Range("A4").Font.Color = 15570276

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

$ws2.Range("A2").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-02-16 10:22:53"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-02-16 10:22:53"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/e2e/41f7d767-a907-42ee-b2c9-5f80e94a7af9.md", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/out/41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/e2e/ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md", "", "", "ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/out/41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f3066fa0cb0df4588b1a40478de660d455eaaa4d/.localization-config", "", "", ".localization-config") | Out-Null

$ws2.Range("A2").Font.Underline = 2
$ws2.Range("A2").Font.Color = 15570276
$ws2.Range("C2").Font.Underline = 2
$ws2.Range("C2").Font.Color = 15570276
$ws2.Range("A3").Font.Underline = 2
$ws2.Range("A3").Font.Color = 15570276
$ws2.Range("C3").Font.Underline = 2
$ws2.Range("C3").Font.Color = 15570276
$ws2.Range("A4").Font.Underline = 2
$ws2.Range("A4").Font.Color = 15570276

# D2/D3/D4 & G2/G3/G4 keep the workbook's date-time number format.
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

$ws3.Range("A2").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf"
$ws3.Range("D2").Value = "2016-02-16 10:23:06"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf"
$ws3.Range("D3").Value = "2016-02-16 10:23:06"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/e2e/41f7d767-a907-42ee-b2c9-5f80e94a7af9.md", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/out/41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/e2e/ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md", "", "", "ffff473fb4a2-2dc3-4dfc-b9dc-5bb4f9a660ba.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTest/oltest/blob/c97fe14d0c784df19f1b84be13b2da20b3a6025c/out/41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf", "", "", "41f7d767-a907-42ee-b2c9-5f80e94a7af9.c97fe14d0c784df19f1b84be13b2da20b3a6025c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f3066fa0cb0df4588b1a40478de660d455eaaa4d/.localization-config", "", "", ".localization-config") | Out-Null

$ws3.Range("A2").Font.Underline = 2
$ws3.Range("A2").Font.Color = 15570276
$ws3.Range("C2").Font.Underline = 2
$ws3.Range("C2").Font.Color = 15570276
$ws3.Range("A3").Font.Underline = 2
$ws3.Range("A3").Font.Color = 15570276
$ws3.Range("C3").Font.Underline = 2
$ws3.Range("C3").Font.Color = 15570276
$ws3.Range("A4").Font.Underline = 2
$ws3.Range("A4").Font.Color = 15570276

$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
